$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 0.3926704299298265
$ws.Range("B3").Value  = 46.55842666749766
$ws.Range("B4").Value  = 786.0167343422662
$ws.Range("B5").Value  = 74.66105371942267
$ws.Range("B6").Value  = 27114.07624177271
$ws.Range("B7").Value  = 717.1778502335162
$ws.Range("B8").Value  = -4516.862650499451
$ws.Range("B9").Value  = 500
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = -4.257500539288378
$ws.Range("B14").Value = -4.197837869044198
$ws.Range("B15").Value = -0.9167842205551707
